# Apply a row permutation to columns D, M, N, O, P, R, S on the active sheet.
# The data in these columns for rows 2-13 is shuffled according to the
# mapping below (target row -> source row, referring to the ORIGINAL data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row number -> old row number (source of the values)
$map = @{
    2  = 8
    3  = 10
    4  = 6
    5  = 4
    6  = 7
    7  = 13
    8  = 2
    9  = 11
    10 = 9
    11 = 3
    12 = 12
    13 = 5
}

# Capture original values for the affected columns before overwriting anything.
$cols = @("D", "M", "N", "O", "P", "R", "S")

$original = @{}
foreach ($row in 2..13) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $original[$row] = $rowData
}

# Write the permuted values back.
foreach ($newRow in $map.Keys) {
    $oldRow = $map[$newRow]
    $rowData = $original[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $rowData[$col]
    }
}
